$wb = $excel.ActiveWorkbook

# Rename sheets to add clarifying suffixes
$wb.Worksheets.Item("Simulation Conditions").Name = "Simulation Conditions (Weather)"
$wb.Worksheets.Item("Rocket Parameters").Name = "Rocket Parameters (Mass)"
$wb.Worksheets.Item("Propellant Parameters").Name = "Propellant Parameters (Tanks)"

# Fix the fuel's chemical formula (Ethane is C2H6, not C2H4) and move the
# active selection, on the renamed "Propellant Parameters (Tanks)" sheet.
$ws = $wb.Worksheets.Item("Propellant Parameters (Tanks)")
$ws.Range("C3").Value = "C2H6"
$ws.Range("H30").Select()
